$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("A12").Value = 130872695
$ws.Range("B12").Value = 79001
$ws.Range("E12").Value = 228912
$ws.Range("F12").Value = "Mörk kolflarnlav"
$ws.Range("G12").Value = "Carbonicola myrmecina"
$ws.Range("H12").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q12").Value = 570816
$ws.Range("R12").Value = 6736802

# Row 13
$ws.Range("A13").Value = 130872698
$ws.Range("B13").Value = 79000
$ws.Range("E13").Value = 6446
$ws.Range("F13").Value = "Kolflarnlav"
$ws.Range("G13").Value = "Carbonicola anthracophila"
$ws.Range("H13").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q13").Value = 570821
$ws.Range("R13").Value = 6736787

# Row 14
$ws.Range("A14").Value = 130872717
$ws.Range("B14").Value = 79243
$ws.Range("E14").Value = 6425
$ws.Range("F14").Value = "Garnlav"
$ws.Range("G14").Value = "Alectoria sarmentosa"
$ws.Range("H14").Value = "(Ach.) Ach."
$ws.Range("Q14").Value = 571254
$ws.Range("R14").Value = 6736578

# Row 25
$ws.Range("A25").Value = 130983063
$ws.Range("B25").Value = 8451
$ws.Range("E25").Value = 106545
$ws.Range("F25").Value = "Mindre märgborre"
$ws.Range("G25").Value = "Tomicus minor"
$ws.Range("H25").Value = "(Hartig, 1834)"
$ws.Range("J25").Value = ""
$ws.Range("K25").Value = ""
$ws.Range("L25").Value = ""
$ws.Range("M25").Value = "äldre gnagspår"
$ws.Range("N25").Value = ""
$ws.Range("Q25").Value = 570956
$ws.Range("R25").Value = 6736657
$ws.Range("S25").Value = 10
$ws.Range("Z25").Value = "09:32"
$ws.Range("AB25").Value = "09:32"
$ws.Range("AC25").Value = ""
$ws.Range("AF25").Value = ""
$ws.Range("AW25").Value = "Bo karlstens"
$ws.Range("AX25").Value = "Bo karlstens"

# Row 26
$ws.Range("A26").Value = 130979083
$ws.Range("B26").Value = 57073
$ws.Range("E26").Value = 100138
$ws.Range("F26").Value = "Tjäder"
$ws.Range("G26").Value = "Tetrao urogallus"
$ws.Range("H26").Value = "Linnaeus, 1758"
$ws.Range("J26").Value = ""
$ws.Range("K26").Value = ""
$ws.Range("L26").Value = ""
$ws.Range("M26").Value = ""
$ws.Range("N26").Value = ""
$ws.Range("Q26").Value = 570745
$ws.Range("R26").Value = 6736794
$ws.Range("S26").Value = 1
$ws.Range("Z26").Value = ""
$ws.Range("AB26").Value = ""
$ws.Range("AC26").Value = "Färsk spillning"
$ws.Range("AF26").Value = ""
$ws.Range("AW26").Value = "Erik Danielsson"
$ws.Range("AX26").Value = "Erik Danielsson"

# Row 32
$ws.Range("A32").Value = 130983060
$ws.Range("B32").Value = 8451
$ws.Range("D32").Value = "LC"
$ws.Range("E32").Value = 106545
$ws.Range("F32").Value = "Mindre märgborre"
$ws.Range("G32").Value = "Tomicus minor"
$ws.Range("H32").Value = "(Hartig, 1834)"
$ws.Range("J32").Value = ""
$ws.Range("K32").Value = ""
$ws.Range("L32").Value = ""
$ws.Range("M32").Value = "äldre gnagspår"
$ws.Range("N32").Value = ""
$ws.Range("P32").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q32").Value = 570988
$ws.Range("R32").Value = 6736721
$ws.Range("Z32").Value = "11:29"
$ws.Range("AB32").Value = "11:29"
$ws.Range("AF32").Value = ""
$ws.Range("AW32").Value = "Bo karlstens"
$ws.Range("AX32").Value = "Bo karlstens"

# Row 33
$ws.Range("A33").Value = 130983618
$ws.Range("B33").Value = 79243
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = "Garnlav"
$ws.Range("G33").Value = "Alectoria sarmentosa"
$ws.Range("H33").Value = "(Ach.) Ach."
$ws.Range("J33").Value = ""
$ws.Range("K33").Value = ""
$ws.Range("L33").Value = ""
$ws.Range("M33").Value = ""
$ws.Range("N33").Value = ""
$ws.Range("P33").Value = "Flytjärnsmyren, Dlr"
$ws.Range("Q33").Value = 570808
$ws.Range("R33").Value = 6736568
$ws.Range("Z33").Value = "09:07"
$ws.Range("AB33").Value = "09:07"
$ws.Range("AF33").Value = ""
$ws.Range("AW33").Value = "Göran Ehn"
$ws.Range("AX33").Value = "Göran Ehn"

# Row 36
$ws.Range("A36").Value = 130979080
$ws.Range("B36").Value = 57884
$ws.Range("E36").Value = 100109
$ws.Range("F36").Value = "Tretåig hackspett"
$ws.Range("G36").Value = "Picoides tridactylus"
$ws.Range("H36").Value = "(Linnaeus, 1758)"
$ws.Range("Q36").Value = 571221
$ws.Range("R36").Value = 6736517
$ws.Range("S36").Value = 1
$ws.Range("Z36").Value = ""
$ws.Range("AB36").Value = ""
$ws.Range("AC36").Value = "Äldre ringhack"
$ws.Range("AF36").Value = ""
$ws.Range("AW36").Value = "Erik Danielsson"
$ws.Range("AX36").Value = "Erik Danielsson"

# Row 37
$ws.Range("A37").Value = 130983068
$ws.Range("B37").Value = 79243
$ws.Range("E37").Value = 6425
$ws.Range("F37").Value = "Garnlav"
$ws.Range("G37").Value = "Alectoria sarmentosa"
$ws.Range("H37").Value = "(Ach.) Ach."
$ws.Range("Q37").Value = 570849
$ws.Range("R37").Value = 6736706
$ws.Range("S37").Value = 10
$ws.Range("Z37").Value = "09:25"
$ws.Range("AB37").Value = "09:25"
$ws.Range("AC37").Value = ""
$ws.Range("AF37").Value = ""
$ws.Range("AW37").Value = "Bo karlstens"
$ws.Range("AX37").Value = "Bo karlstens"

# Row 38
$ws.Range("A38").Value = 130983600
$ws.Range("B38").Value = 57884
$ws.Range("E38").Value = 100109
$ws.Range("F38").Value = "Tretåig hackspett"
$ws.Range("G38").Value = "Picoides tridactylus"
$ws.Range("H38").Value = "(Linnaeus, 1758)"
$ws.Range("K38").Value = ""
$ws.Range("L38").Value = ""
$ws.Range("M38").Value = "äldre spår"
$ws.Range("N38").Value = ""
$ws.Range("P38").Value = "Flytjärnsmyren, Dlr"
$ws.Range("Q38").Value = 571022
$ws.Range("R38").Value = 6736648
$ws.Range("Z38").Value = "09:33"
$ws.Range("AB38").Value = "09:33"
$ws.Range("AE38").Value = $true
$ws.Range("AF38").Value = ""
$ws.Range("AW38").Value = "Göran Ehn"
$ws.Range("AX38").Value = "Göran Ehn"

# Row 39
$ws.Range("A39").Value = 130983070
$ws.Range("P39").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q39").Value = 570811
$ws.Range("R39").Value = 6736543
$ws.Range("Z39").Value = "09:04"
$ws.Range("AB39").Value = "09:04"
$ws.Range("AF39").Value = ""
$ws.Range("AW39").Value = "Bo karlstens"
$ws.Range("AX39").Value = "Bo karlstens"

# Row 40
$ws.Range("A40").Value = 130983617
$ws.Range("B40").Value = 79243
$ws.Range("E40").Value = 6425
$ws.Range("F40").Value = "Garnlav"
$ws.Range("G40").Value = "Alectoria sarmentosa"
$ws.Range("H40").Value = "(Ach.) Ach."
$ws.Range("K40").Value = ""
$ws.Range("L40").Value = ""
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = ""
$ws.Range("Q40").Value = 570789
$ws.Range("R40").Value = 6736672
$ws.Range("Z40").Value = "09:12"
$ws.Range("AB40").Value = "09:12"
$ws.Range("AE40").Value = $false
